$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case row (row 2): ID, host, interface, method, header(blank), datas, statue_code, except(blank)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "http://192.168.1.86:80"
[void]$ws.Hyperlinks.Add($ws.Range("B2"), "http://192.168.1.86:80")
$ws.Range("C2").Value = "/qc-engine/v3/sysusers/login"
$ws.Range("D2").Value = "POST"
$ws.Range("F2").Value = '{"userName": "admin", "password": "123"}'
$ws.Range("G2").Value = 200
$ws.Range("H2").Value = 200

# New "msg" header column
$ws.Range("I1").Value = "msg"

# Widen host column so the URL is readable
$ws.Columns.Item(2).ColumnWidth = 23.43

# Leave selection where the author ended up
[void]$ws.Range("J3").Select()
